# Edit script: fix descriptives example data set for pooled variables.
# 1) Rename the existing ordinal "pvkat_*" block (rows 20-49, column A) to "pvord_*"
#    (these rows keep their existing 5-category structure: Kompetenzstufe 1-5).
# 2) Append a brand-new categorical/nominal "pvkat_*" block (rows 50-85) that
#    mirrors the same structure but with an extra 6th category (Kompetenzstufe 6)
#    so the example data covers a pooled variable with differing per-imputation
#    category counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: rename column A entries for rows 20-49 from pvkat_* to pvord_* ---
for ($r = 20; $r -le 49; $r++) {
    $cur = $ws.Cells.Item($r, 1).Value2
    $new = $cur -replace '^pvkat_', 'pvord_'
    $ws.Cells.Item($r, 1).Value = $new
}

# --- Step 2: append the new pvkat_* rows (50-85) ---
$newRows = @(
    @("pvkat_1", 1, "nein", "Kompetenzstufe 1", "nein"),
    @("pvkat_1", 2, "nein", "Kompetenzstufe 2", "nein"),
    @("pvkat_1", 3, "nein", "Kompetenzstufe 3", "nein"),
    @("pvkat_1", 4, "nein", "Kompetenzstufe 4", "nein"),
    @("pvkat_1", 5, "nein", "Kompetenzstufe 5", "nein"),
    @("pvkat_pooled", 1, "nein", "Kompetenzstufe 1", "nein"),
    @("pvkat_pooled", 2, "nein", "Kompetenzstufe 2", "nein"),
    @("pvkat_pooled", 3, "nein", "Kompetenzstufe 3", "nein"),
    @("pvkat_pooled", 4, "nein", "Kompetenzstufe 4", "nein"),
    @("pvkat_pooled", 5, "nein", "Kompetenzstufe 5", "nein"),
    @("pvkat_pooled", 6, "nein", "Kompetenzstufe 6", "nein"),
    @("pvkat_1", 6, "nein", "Kompetenzstufe 6", "nein"),
    @("pvkat_2", 1, "nein", "Kompetenzstufe 1", "nein"),
    @("pvkat_2", 2, "nein", "Kompetenzstufe 2", "nein"),
    @("pvkat_2", 3, "nein", "Kompetenzstufe 3", "nein"),
    @("pvkat_2", 4, "nein", "Kompetenzstufe 4", "nein"),
    @("pvkat_2", 5, "nein", "Kompetenzstufe 5", "nein"),
    @("pvkat_2", 6, "nein", "Kompetenzstufe 6", "nein"),
    @("pvkat_3", 1, "nein", "Kompetenzstufe 1", "nein"),
    @("pvkat_3", 2, "nein", "Kompetenzstufe 2", "nein"),
    @("pvkat_3", 3, "nein", "Kompetenzstufe 3", "nein"),
    @("pvkat_3", 4, "nein", "Kompetenzstufe 4", "nein"),
    @("pvkat_3", 5, "nein", "Kompetenzstufe 5", "nein"),
    @("pvkat_3", 6, "nein", "Kompetenzstufe 6", "nein"),
    @("pvkat_4", 1, "nein", "Kompetenzstufe 1", "nein"),
    @("pvkat_4", 2, "nein", "Kompetenzstufe 2", "nein"),
    @("pvkat_4", 3, "nein", "Kompetenzstufe 3", "nein"),
    @("pvkat_4", 4, "nein", "Kompetenzstufe 4", "nein"),
    @("pvkat_4", 5, "nein", "Kompetenzstufe 5", "nein"),
    @("pvkat_4", 6, "nein", "Kompetenzstufe 6", "nein"),
    @("pvkat_5", 1, "nein", "Kompetenzstufe 1", "nein"),
    @("pvkat_5", 2, "nein", "Kompetenzstufe 2", "nein"),
    @("pvkat_5", 3, "nein", "Kompetenzstufe 3", "nein"),
    @("pvkat_5", 4, "nein", "Kompetenzstufe 4", "nein"),
    @("pvkat_5", 5, "nein", "Kompetenzstufe 5", "nein"),
    @("pvkat_5", 6, "nein", "Kompetenzstufe 6", "nein")
)

$r = 50
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}
